# Daily Report update: append new daily snapshot (date serial 46070) to
# Daily_Data, then refresh the Today_Summary and Monthly_Stats roll-ups
# to reflect it.

$wb  = $excel.ActiveWorkbook
$wsDaily   = $wb.Worksheets.Item("Daily_Data")
$wsToday   = $wb.Worksheets.Item("Today_Summary")
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")

# ---------------------------------------------------------------------
# 1) Daily_Data: append rows 146-169 for date 46070
#    columns: A=Date B=Region_Type C=PREV_TOTAL D=RECEIVED E=WITHDRAWN
#             F=NET_CHANGE G=ADJUSTMENT H=TOTAL_TODAY
# ---------------------------------------------------------------------

$newDate = 46070

$dailyRows = @(
    @("ASAHI DEPOSITORY LLC Registered", 23953631.592, 0, 0, 0, -586856.4, 23366775.192),
    @("ASAHI DEPOSITORY LLC Eligible", 2097038.208, 0, 0, 0, 586856.4, 2683894.608),
    @("BRINK'S, INC. Registered", 16122359.646, 0, 0, 0, -344717.66, 15777641.986),
    @("BRINK'S, INC. Eligible", 39587772.794, 0, 0, 0, 344717.66, 39932490.454),
    @("CNT DEPOSITORY, INC. Registered", 12235256.378, 0, 0, 0, -60404.809, 12174851.569),
    @("CNT DEPOSITORY, INC. Eligible", 15027218.389, 0, 912481.262, -912481.262, 60404.809, 14175141.936),
    @("DELAWARE DEPOSITORY Registered", 1547695.233, 0, 0, 0, -14918.81, 1532776.423),
    @("DELAWARE DEPOSITORY Eligible", 16254567.062, 0, 0, 0, 14918.81, 16269485.872),
    @("HSBC BANK, USA Registered", 3472271.68, 0, 0, 0, -60114.11, 3412157.57),
    @("HSBC BANK, USA Eligible", 21150312.483, 0, 305487.68, -305487.68, 60114.11, 20904938.913),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 273789.87, 0, 0, 0, 0, 273789.87),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 3642206.244, 0, 0, 0, 0, 3642206.244),
    @("JP MORGAN CHASE BANK NA Registered", 12025248.54, 0, 0, 0, -24904.77, 12000343.77),
    @("JP MORGAN CHASE BANK NA Eligible", 145773316.663, 0, 1943813, -1943813, 24904.77, 143854408.433),
    @("LOOMIS INTERNATIONAL (US) LLC Registered", 7393354.417, 0, 0, 0, -540734.24, 6852620.177),
    @("LOOMIS INTERNATIONAL (US) LLC Eligible", 22952116.706, 0, 0, 0, 540734.24, 23492850.946),
    @("MALCA-AMIT ARMORED, INC. Registered", 0, 0, 0, 0, 0, 0),
    @("MALCA-AMIT ARMORED, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @("MALCA-AMIT USA, LLC Registered", 1225506.264, 0, 0, 0, -275872.2, 949634.064),
    @("MALCA-AMIT USA, LLC Eligible", 798026.177, 0, 0, 0, 275872.2, 1073898.377),
    @("MANFRA, TORDELLA & BROOKES, LLC Registered", 6365432.292, 0, 0, 0, -145802.259, 6219630.033),
    @("MANFRA, TORDELLA & BROOKES, LLC Eligible", 12302849.048, 0, 0, 0, 145802.259, 12448651.307),
    @("STONEX PRECIOUS METALS LLC Registered", 7540323.54, 0, 0, 0, -1308822.14, 6231501.4),
    @("STONEX PRECIOUS METALS LLC Eligible", 233197.18, 0, 0, 0, 1308822.14, 1542019.32)
)

$dateFmt = $wsDaily.Cells.Item(145, 1).NumberFormat

$startRow = 146
$r = $startRow
foreach ($row in $dailyRows) {
    $wsDaily.Cells.Item($r, 1).Value = $newDate
    $wsDaily.Cells.Item($r, 1).NumberFormat = $dateFmt
    $wsDaily.Cells.Item($r, 2).Value = $row[0]
    $wsDaily.Cells.Item($r, 3).Value = $row[1]
    $wsDaily.Cells.Item($r, 4).Value = $row[2]
    $wsDaily.Cells.Item($r, 5).Value = $row[3]
    $wsDaily.Cells.Item($r, 6).Value = $row[4]
    $wsDaily.Cells.Item($r, 7).Value = $row[5]
    $wsDaily.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Build a lookup of Region_Type -> TOTAL_TODAY (H) and -> WITHDRAWN (E) for
# the newly appended date, used to refresh the summary sheets below.
$totalToday = @{}
$withdrawnToday = @{}
foreach ($row in $dailyRows) {
    $totalToday[$row[0]] = $row[6]
    $withdrawnToday[$row[0]] = $row[3]
}

# ---------------------------------------------------------------------
# 2) Today_Summary: refresh Eligible(B) / Registered(C) / Total_Stock(D)
#    for each depository using the freshly-appended TOTAL_TODAY values.
# ---------------------------------------------------------------------

$depositories = @(
    "ASAHI DEPOSITORY LLC",
    "BRINK'S, INC.",
    "CNT DEPOSITORY, INC.",
    "DELAWARE DEPOSITORY",
    "HSBC BANK, USA",
    "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE",
    "JP MORGAN CHASE BANK NA",
    "LOOMIS INTERNATIONAL (US) LLC",
    "MALCA-AMIT ARMORED, INC.",
    "MALCA-AMIT USA, LLC",
    "MANFRA, TORDELLA & BROOKES, LLC",
    "STONEX PRECIOUS METALS LLC"
)

$sumRow = 2
foreach ($dep in $depositories) {
    $eligible = $totalToday[$dep + " Eligible"]
    $registered = $totalToday[$dep + " Registered"]
    $wsToday.Cells.Item($sumRow, 2).Value = $eligible
    $wsToday.Cells.Item($sumRow, 3).Value = $registered
    $wsToday.Cells.Item($sumRow, 4).Value = $eligible + $registered
    $sumRow = $sumRow + 1
}

# ---------------------------------------------------------------------
# 3) Monthly_Stats: refresh the per-depository detail block (rows 7-30)
#    and the month grand-total block (row 2).
# ---------------------------------------------------------------------

$monthlyNames = @(
    "ASAHI DEPOSITORY LLC Eligible",
    "ASAHI DEPOSITORY LLC Registered",
    "BRINK'S, INC. Eligible",
    "BRINK'S, INC. Registered",
    "CNT DEPOSITORY, INC. Eligible",
    "CNT DEPOSITORY, INC. Registered",
    "DELAWARE DEPOSITORY Eligible",
    "DELAWARE DEPOSITORY Registered",
    "HSBC BANK, USA Eligible",
    "HSBC BANK, USA Registered",
    "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible",
    "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered",
    "JP MORGAN CHASE BANK NA Eligible",
    "JP MORGAN CHASE BANK NA Registered",
    "LOOMIS INTERNATIONAL (US) LLC Eligible",
    "LOOMIS INTERNATIONAL (US) LLC Registered",
    "MALCA-AMIT ARMORED, INC. Eligible",
    "MALCA-AMIT ARMORED, INC. Registered",
    "MALCA-AMIT USA, LLC Eligible",
    "MALCA-AMIT USA, LLC Registered",
    "MANFRA, TORDELLA & BROOKES, LLC Eligible",
    "MANFRA, TORDELLA & BROOKES, LLC Registered",
    "STONEX PRECIOUS METALS LLC Eligible",
    "STONEX PRECIOUS METALS LLC Registered"
)

$eligibleTotal = 0
$registeredTotal = 0

$monthRow = 7
foreach ($name in $monthlyNames) {
    $withdrawnDelta = $withdrawnToday[$name]
    if ($withdrawnDelta -ne 0) {
        $oldWithdrawn = $wsMonthly.Cells.Item($monthRow, 4).Value2
        $wsMonthly.Cells.Item($monthRow, 4).Value = $oldWithdrawn + $withdrawnDelta
    }
    $newTotal = $totalToday[$name]
    $wsMonthly.Cells.Item($monthRow, 5).Value = $newTotal

    if ($name -like "* Eligible") {
        $eligibleTotal = $eligibleTotal + $newTotal
    } else {
        $registeredTotal = $registeredTotal + $newTotal
    }

    $monthRow = $monthRow + 1
}

$wsMonthly.Cells.Item(2, 2).Value = $eligibleTotal
$wsMonthly.Cells.Item(2, 3).Value = $registeredTotal
$wsMonthly.Cells.Item(2, 4).Value = $eligibleTotal + $registeredTotal
